# Atualização automática de preços de eletricidade
# Updates row 2 of the spot price worksheet with the latest daily values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46008
$ws.Range("B2").Value = 83.14
$ws.Range("C2").Value = 82.26000000000001
$ws.Range("D2").Value = 81.20999999999999
$ws.Range("E2").Value = 77.53
$ws.Range("F2").Value = 76.53
$ws.Range("G2").Value = 78.59
$ws.Range("H2").Value = 97.93000000000001
$ws.Range("I2").Value = 113.82
$ws.Range("J2").Value = 136.16
$ws.Range("K2").Value = 114.85
$ws.Range("L2").Value = 98.45999999999999
$ws.Range("M2").Value = 90.67
$ws.Range("N2").Value = 87.86
$ws.Range("O2").Value = 84.38
$ws.Range("P2").Value = 86.23999999999999
$ws.Range("Q2").Value = 94.01000000000001
$ws.Range("R2").Value = 110.13
$ws.Range("S2").Value = 112.4
$ws.Range("T2").Value = 124.22
$ws.Range("U2").Value = 123.11
$ws.Range("V2").Value = 136.16
$ws.Range("W2").Value = 123.48
$ws.Range("X2").Value = 115.14
$ws.Range("Y2").Value = 105.66
$ws.Range("Z2").Value = 101.41
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 120.11
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 129.82
$ws.Range("AE2").Value = "8h-10h"
$ws.Range("AF2").Value = 125.5
$ws.Range("AG2").Value = "0h-15h"
